# Extend the Gradient_Month_12 data range from A1:G20 to A1:G25 by appending
# five new distance rows (1100, 1200, 1300, 1400, 1500) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data: distance, MEAN, STD, MIN, MAX, COUNT, Month
$newRows = @(
    @(1100, 4.948644638061523, 1.441019535064697, -2.955584764480591, 12.23557758331299, 18372, "12"),
    @(1200, 4.90306282043457,  1.523025751113892, -3.123027563095093, 11.88702392578125, 18502, "12"),
    @(1300, 4.851481914520264, 1.523171305656433, -6.16433572769165,  11.98783111572266, 18317, "12"),
    @(1400, 4.851562023162842, 1.553480625152588, -6.273685932159424, 10.66537475585938, 18231, "12"),
    @(1500, 4.857754707336426, 1.714037775993347, -6.256599903106689, 15.21708488464355, 18189, "12")
)

$startRow = 21
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]

    # Column G ("Month") is stored as text ("12"), not a number, in the
    # source data - force text formatting, assign, then restore the
    # default "Normal" style so no stray numFmt/style index lingers on
    # the cell (matches the plain, style-less cells used elsewhere).
    $gCell = $ws.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $rowData[6]
    $gCell.Style = "Normal"
}
